$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (Enterprise Project) below the existing rows
$ws.Range("A6").Value = "Enterprise Proj sheet"
$ws.Range("B6").Value = "1.1EnterpriseProject"

# Copy formatting from row 5 (Secgroup sheet) to keep consistent styling for the new row
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to match the new state of the sheet
$ws.Range("B8").Select()
